# Add data for 2021-11-28
# - Rename the sheet/title from "Through November 19" to "Through November 20"
# - Update the corresponding column header text
# - Update counts for the current month's column (and a few other historical
#   cells, per the upstream diff) for several neighborhoods

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename worksheet tab to reflect new "through" date
$ws.Name = "Through 2021-11-20"

# Update the header text for column B (the "current month" column)
$ws.Range("B1").Value = "November 2021 (through November 20)"

# --- Cell value updates (row => neighborhood, column => month) ---

# North Lawndale (row 2)
$ws.Range("B2").Value = 9
$ws.Range("M2").Value = 12

# Garfield Park (row 3)
$ws.Range("M3").Value = 11

# Humboldt Park (row 4)
$ws.Range("BE4").Value = 3

# South Shore (row 6)
$ws.Range("B6").Value = 8

# Englewood (row 8)
$ws.Range("AI8").Value = 2
$ws.Range("BE8").Value = 4

# Auburn Gresham (row 9)
$ws.Range("M9").Value = 5

# Grand Boulevard (row 12)
$ws.Range("M12").Value = 6

# Row 16
$ws.Range("BP16").Value = 1

# Row 17 (West Loop)
$ws.Range("AI17").Value = 1

# Row 20 (Chatham)
$ws.Range("M20").Value = 2
$ws.Range("AT20").Value = 1
$ws.Range("BP20").Value = 1

# Row 29 (Near South Side)
$ws.Range("AT29").Value = 1

# Row 34 (Hyde Park)
$ws.Range("AT34").Value = 2

# Row 37 (Ashburn)
$ws.Range("BE37").Value = 4

# Row 39 (Avalon Park)
$ws.Range("M39").Value = 1

# Row 48 (Roseland)
$ws.Range("B48").Value = 5
$ws.Range("M48").Value = 4

# Row 49 (Logan Square)
$ws.Range("AT49").Value = 1

# Row 52 (Oakland)
$ws.Range("B52").Value = 1

# Row 55 (North Center)
$ws.Range("M55").Value = 1

# Row 74 (Garfield Ridge)
$ws.Range("B74").Value = 1

# Row 89 (Portage Park)
$ws.Range("AI89").Value = 2

# Row 98 (Woodlawn)
$ws.Range("B98").Value = 2
